# "Mise à jour de l'application" — record a new CDF match ("CDF T6") in the
# player log. The header columns HK:HN already exist as a 4-column match
# block (match-minutes / T-R-NR-HG status / goals / assists); HK1 currently
# just says "CDF" and is renamed to the specific "CDF T6" label. Everything
# else in this block (the attendance status for every player, plus minutes
# / goals / assists where relevant) is new data for that match. Two players
# also picked up a previously-missing "R2J5" (Réserve) appearance at JK/JL.
#
# All the summary / COUNTIF / SUM formulas already span these columns by
# range, so they recalculate automatically once the underlying data is
# written — no formula cells need to be touched directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the match header (was generic "CDF", now "CDF T6")
$ws.Range("HK1").Value = "CDF T6"

# Row 2 - Alban Rambaud: hors-groupe -> did not enter (NR) for CDF T6
$ws.Range("HL2").Value = "NR"

# Row 3 - Jassim Assoul: HG for CDF T6; also a newly logged R2J5 (90 min, T)
$ws.Range("HL3").Value = "HG"
$ws.Range("JK3").Value = 90
$ws.Range("JL3").Value = "T"

# Row 4 - Enzo Vita: started CDF T6, played 90 min
$ws.Range("HK4").Value = 90
$ws.Range("HL4").Value = "T"

# Row 5 - Romain Thunet: started CDF T6, played 90 min
$ws.Range("HK5").Value = 90
$ws.Range("HL5").Value = "T"

# Row 6 - Amine Taiar: HG for CDF T6
$ws.Range("HL6").Value = "HG"

# Row 7 - Naim Ighbane: started CDF T6, played 90 min
$ws.Range("HK7").Value = 90
$ws.Range("HL7").Value = "T"

# Row 8 - Hedi Nasri: HG for CDF T6
$ws.Range("HL8").Value = "HG"

# Row 9 - Mattheo Haon: entered as sub (R) for CDF T6, played 26 min, 1 assist
$ws.Range("HK9").Value = 26
$ws.Range("HL9").Value = "R"
$ws.Range("HN9").Value = 1

# Row 10 - Maé Clavel: entered as sub (R) for CDF T6, played 45 min, 1 goal
$ws.Range("HK10").Value = 45
$ws.Range("HL10").Value = "R"
$ws.Range("HM10").Value = 1

# Row 11 - Levy Ndoutoume: started CDF T6, played 64 min
$ws.Range("HK11").Value = 64
$ws.Range("HL11").Value = "T"

# Row 12 - Yanis Berrached: HG for CDF T6
$ws.Range("HL12").Value = "HG"

# Row 13 - Rayane Chayebi: HG for CDF T6
$ws.Range("HL13").Value = "HG"

# Row 14 - Ilan Ihaddadene: started CDF T6, played 90 min, 1 assist
$ws.Range("HK14").Value = 90
$ws.Range("HL14").Value = "T"
$ws.Range("HN14").Value = 1

# Row 15 - Karahali Souaré: entered as sub (R) for CDF T6, played 26 min, 1 assist
$ws.Range("HK15").Value = 26
$ws.Range("HL15").Value = "R"
$ws.Range("HN15").Value = 1

# Row 16 - Amir Etien: entered as sub (R) for CDF T6, played 26 min, 1 goal
$ws.Range("HK16").Value = 26
$ws.Range("HL16").Value = "R"
$ws.Range("HM16").Value = 1

# Row 17 - Karim Belmahi: HG for CDF T6
$ws.Range("HL17").Value = "HG"

# Row 18 - Emmanuel Valey: started CDF T6, played 45 min, 1 assist
$ws.Range("HK18").Value = 45
$ws.Range("HL18").Value = "T"
$ws.Range("HN18").Value = 1

# Row 19 - Jeremie Laurent: HG for CDF T6
$ws.Range("HL19").Value = "HG"

# Row 20 - Sofiane Belle: started CDF T6, played 90 min, 4 goals, 1 assist
$ws.Range("HK20").Value = 90
$ws.Range("HL20").Value = "T"
$ws.Range("HM20").Value = 4
$ws.Range("HN20").Value = 1

# Row 21 - Amir Kherrab: HG for CDF T6
$ws.Range("HL21").Value = "HG"

# Row 22 - Naim Dhib: started CDF T6, played 64 min, 1 goal, 1 assist
$ws.Range("HK22").Value = 64
$ws.Range("HL22").Value = "T"
$ws.Range("HM22").Value = 1
$ws.Range("HN22").Value = 1

# Row 23 - Wael Fareh: HG for CDF T6
$ws.Range("HL23").Value = "HG"

# Row 24 - Yoan Zouma: HG for CDF T6
$ws.Range("HL24").Value = "HG"

# Row 25 - Ilyes Boughanmi: started CDF T6, played 90 min, 3 goals
$ws.Range("HK25").Value = 90
$ws.Range("HL25").Value = "T"
$ws.Range("HM25").Value = 3

# Row 26 - Omar Benyounes: HG for CDF T6; also a newly logged R2J5 (90 min, T)
$ws.Range("HL26").Value = "HG"
$ws.Range("JK26").Value = 90
$ws.Range("JL26").Value = "T"

# Row 27 - Yoann Martelat: started CDF T6, played 90 min
$ws.Range("HK27").Value = 90
$ws.Range("HL27").Value = "T"

# Row 28 - Malik Boussaid: started CDF T6, played 64 min, 2 assists
$ws.Range("HK28").Value = 64
$ws.Range("HL28").Value = "T"
$ws.Range("HN28").Value = 2

# Row 29 - Kamal Bafounta: HG for CDF T6
$ws.Range("HL29").Value = "HG"

# Recalculate all the SUM/COUNTIF summary columns now that the raw match
# data is in place.
$excel.CalculateFullRebuild()

# Leave the selection near the newly-entered data, matching where the
# author ended up after this edit.
$ws.Range("JP22").Select()
